$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 155
$ws.Range("I2").Value = 437
$ws.Range("J2").Value = 1830
$ws.Range("K2").Value = 7
$ws.Range("L2").Value = 504
$ws.Range("M2").Value = 33
$ws.Range("N2").Value = 303
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 6
$ws.Range("R2").Value = 29
$ws.Range("S2").Value = 177
$ws.Range("T2").Value = 306
$ws.Range("U2").Value = 26
$ws.Range("V2").Value = 2843
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 2758
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 47
$ws.Range("AA2").Value = 26
